$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Header row (row 1): update wording of existing headers and add a new
# column G header for the newly introduced "سازه" (structure) gauge space ---
$ws2.Cells.Item(1, 3).Value = "قابلیت عبور از فضای مجاز"
$ws2.Cells.Item(1, 5).Value = "قابلیت عبور از فضای آزاد"
$ws2.Cells.Item(1, 6).Value = "اندازه ورود به فضای سازه"

# Create the G1 header cell with the same look as the rest of the header row
$ws2.Range("F1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$ws2.Cells.Item(1, 7).Value = "قابلیت عبور از فضای سازه"

# Give the new column a normal width like the others
$ws2.Columns.Item(7).ColumnWidth = 20.166666666666668

# Row 2 no longer carries a value in column F (its meaning moved to the new
# "سازه" columns used starting on the final destination row)
$ws2.Range("F2").Clear()

# --- New route-segment rows (3-5), styled like row 2 ---
$ws2.Range("A2:E2").Copy()
$ws2.Range("A3:E5").PasteSpecial(-4122)

$ws2.Cells.Item(3, 1).Value = "بهرام - باغ یک"
$ws2.Cells.Item(3, 2).Value = "گاباری 5.4"
$ws2.Cells.Item(3, 3).Value = "غیر قابل عبور"
$ws2.Cells.Item(3, 4).Value = 25.0
$ws2.Cells.Item(3, 5).Value = "قابل عبور"

$ws2.Cells.Item(4, 1).Value = "باغ یک - سواریان"
$ws2.Cells.Item(4, 2).Value = "گاباری 5.2"
$ws2.Cells.Item(4, 3).Value = "غیر قابل عبور"
$ws2.Cells.Item(4, 4).Value = 35.35533905932738
$ws2.Cells.Item(4, 5).Value = "قابل عبور"

$ws2.Cells.Item(5, 1).Value = "سواریان - نورآباد"
$ws2.Cells.Item(5, 2).Value = "گاباری 5.4"
$ws2.Cells.Item(5, 3).Value = "غیر قابل عبور"
$ws2.Cells.Item(5, 4).Value = 25.0
$ws2.Cells.Item(5, 5).Value = "قابل عبور"

# --- Final destination row (6): also fills columns F and G ---
$ws2.Range("A2:E2").Copy()
$ws2.Range("A6:E6").PasteSpecial(-4122)
$ws2.Range("C2").Copy()
$ws2.Range("F6:G6").PasteSpecial(-4122)

$ws2.Cells.Item(6, 1).Value = "نورآباد - بندرامام خمینی"
$ws2.Cells.Item(6, 2).Value = "گاباری ۴.۷"
$ws2.Cells.Item(6, 3).Value = "غیر قابل عبور"
$ws2.Cells.Item(6, 4).Value = 335.4101966249685
$ws2.Cells.Item(6, 5).Value = "غیر قابل عبور"
$ws2.Cells.Item(6, 6).Value = 100.0
$ws2.Cells.Item(6, 7).Value = "غیر قابل عبور"
